$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    3  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    4  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    5  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    6  = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 7.143138311642302)
    7  = @(0.6606524410359556, 0.04071648406533734, 3.537761648806719, 0.4942365360607697, 4.733367109968782)
    8  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    9  = @(1.455362044514542, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 3.009163075608874)
    10 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    11 = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 7.143138311642302)
    12 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
